# Jenkins: Auto-update Test Reports and Screenshots
# Append a new generated user row to the "CreatedUsers" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CreatedUsers")

# Find the next empty row in column A (data starts at row 1)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

$ws.Cells.Item($newRow, 1).Value = "user_1768216212073"
$ws.Cells.Item($newRow, 2).Value = "password123"
